# Weekly update: insert a new daily price record (row 494) for
# "Vega Modelo de Temuco - Acelga" and push the existing historical
# rows (494-510) down by one (495-511).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 494, shifting rows
# 494:510 down to 495:511 (and growing the sheet dimension to R511).
$ws.Rows.Item(494).Insert()

# Populate the newly inserted row 494 with the new record.
$ws.Cells.Item(494, 1).Value = 10
$ws.Cells.Item(494, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(494, 3).Value = "La Araucanía"
$ws.Cells.Item(494, 4).Value = 45075
$ws.Cells.Item(494, 5).Value = 9
$ws.Cells.Item(494, 6).Value = 100112009
$ws.Cells.Item(494, 7).Value = "Acelga"
$ws.Cells.Item(494, 8).Value = "Sin especificar"
$ws.Cells.Item(494, 9).Value = "Primera"
$ws.Cells.Item(494, 10).Value = 110
$ws.Cells.Item(494, 11).Value = 6000
$ws.Cells.Item(494, 12).Value = 6000
$ws.Cells.Item(494, 13).Value = 6000
$ws.Cells.Item(494, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(494, 15).Value = "Región del Maule"
$ws.Cells.Item(494, 16).Value = 500
$ws.Cells.Item(494, 17).Value = 12
$ws.Cells.Item(494, 18).Value = "Hortaliza"
